$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Sapna R.%Kudchadkar%NULL%1,                  Christopher L.%Carroll%NULL%1]'
$ws.Range("E3").Value = '[Sohaib R%Rufai%sohaibrufai@gmail.com%1,                  Catey%Bunce%Catey.bunce@kcl.ac.uk%2,                  Catey%Bunce%Catey.bunce@kcl.ac.uk%0]'
$ws.Range("E4").Value = '[Travis%Sanchez%NULL%2,                  Kamal%Al Nasr%NULL%2,                  Kamal%Al Nasr%NULL%0,                  Ross%Gore%NULL%1,                  Abdullah%Wahbeh%NULL%2,                  Abdullah%Wahbeh%NULL%0,                  Tareq%Nasralah%t.nasralah@northeastern.edu%2,                  Tareq%Nasralah%t.nasralah@northeastern.edu%0,                  Mohammad%Al-Ramahi%NULL%2,                  Mohammad%Al-Ramahi%NULL%0,                  Omar%El-Gayar%NULL%2,                  Omar%El-Gayar%NULL%0]'
$ws.Range("E8").Value = '[Amy%Chesser%xref no email%1,           Amy%Drassen Ham%xref no email%1,           Nikki%Keene Woods%xref no email%1]'
$ws.Range("E9").Value = '[Vittorio%Gebbia%NULL%1,                  Dario%Piazza%NULL%1,                  Maria Rosaria%Valerio%NULL%1,                  Nicolò%Borsellino%NULL%1,                  Alberto%Firenze%NULL%1]'
$ws.Range("E10").Value = '[Jesus-Angel%P\u00e9rez-Dasilva%xref no email%1,           Koldobika%Meso-Ayerdi%xref no email%1,           Terese%Mendiguren-Galdosp\u00edn%xref no email%1]'
$ws.Range("E11").Value = '[Maria Renee%Jimenez‐Sotomayor%NULL%1,                  Carolina%Gomez‐Moreno%NULL%1,                  Enrique%Soto‐Perez‐de‐Celis%enrique.sotop@incmnsz.mx%1]'
$ws.Range("E12").Value = '[Greg%Kawchuk%greg.kawchuk@ualberta.ca%1,                  Jan%Hartvigsen%jhartvigsen@health.sdu.dk%2,                  Jan%Hartvigsen%jhartvigsen@health.sdu.dk%0,                  Stan%Innes%S.Innes@murdoch.edu.au%1,                  J. Keith%Simpson%k.simpson@murdoch.edu.au%1,                  Brian%Gushaty%bgushaty@gushaty.com%1]'
$ws.Range("E13").Value = '[Gunther%Eysenbach%NULL%0,                  Jon-Patrick%Allem%NULL%2,                  Jon-Patrick%Allem%NULL%0,                  Richard%Zowalla%NULL%1,                  Wasim%Ahmed%Wasim.Ahmed@Newcastle.ac.uk%2,                  Wasim%Ahmed%Wasim.Ahmed@Newcastle.ac.uk%0,                  Josep%Vidal-Alaball%NULL%2,                  Josep%Vidal-Alaball%NULL%0,                  Joseph%Downing%NULL%2,                  Joseph%Downing%NULL%0,                  Francesc%López Seguí%NULL%2,                  Francesc%López Seguí%NULL%0]'
$ws.Range("E15").Value = '[Francisco%Morinha%xref no email%1,           Paula%Magalh\u00e3es%xref no email%1]'
$ws.Range("E16").Value = '[Servet%Aker%servetaker@gmail.com%1,                  Özlem%Mıdık%NULL%2,                  Özlem%Mıdık%NULL%0]'
$ws.Range("E17").Value = '[Cristina M%Pulido%xref no email%1,           Beatriz%Villarejo-Carballido%xref no email%1,           Gisela%Redondo-Sama%xref no email%1,           Aitor%G\u00f3mez%xref no email%1]'
$ws.Range("E18").Value = '[Flecha%Ram\u00f3n%coreGivesNoEmail%1,                 Guo%Mengna%coreGivesNoEmail%1,                 Pulido%Rodr\u00edguez Cristina%coreGivesNoEmail%1,                 Ramis%Salas Maria del Mar%coreGivesNoEmail%1,                 Redondo%Sama Gisela%coreGivesNoEmail%1,                 Villarejo%Carballido Beatriz%coreGivesNoEmail%1]'
$ws.Range("E19").Value = '[Bahar%Yuksel%baharyl86@gmail.com%1,                  Kubra%Cakmak%NULL%1]'
$ws.Range("E20").Value = '[Gunther%Eysenbach%NULL%0,                  Nazakat%Hamassed%NULL%2,                  Nazakat%Hamassed%NULL%0,                  Hardawan%Kakashekh%NULL%1,                  Muhammad%Saud%NULL%1,                  Mohammad Amin%Bahrami%NULL%1,                  Araz Ramazan%Ahmad%araz.ahmad85@uor.edu.krd%2,                  Araz Ramazan%Ahmad%araz.ahmad85@uor.edu.krd%0,                  Hersh Rasool%Murad%NULL%2,                  Hersh Rasool%Murad%NULL%0]'
$ws.Range("E21").Value = '[Alexander%Muacevic%NULL%0,                  John R%Adler%NULL%0,                  Ramez%Kouzy%NULL%2,                  Ramez%Kouzy%NULL%0,                  Joseph%Abi Jaoude%NULL%1,                  Afif%Kraitem%NULL%1,                  Molly B%El Alam%NULL%1,                  Basil%Karam%NULL%1,                  Elio%Adib%NULL%1,                  Jabra%Zarka%NULL%1,                  Cindy%Traboulsi%NULL%1,                  Elie W%Akl%NULL%1,                  Khalil%Baddour%NULL%1]'
